$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update example filelist.xlsx contents:
#  - Processor renamed from "Z Vanderbosch" to "J Doe"
#  - Row 2 (ROMO) location label fixed ("Rocky Mountain_NP" -> "Rocky Mountain NP")
#    and its Central_AZ value reset to 0
#  - Row 3 (NIOB) B_band flag flipped to "No"
#  - Row 4's example dataset swapped from the "Zachs Backyard" entry
#    (ZABA250503) to a new "Grand Canyon NP" entry (GRCA120730)
$ws.Range("I2").Value = "J Doe"
$ws.Range("K2").Value = "Rocky Mountain NP"
$ws.Range("D3").Value = "No"
$ws.Range("I3").Value = "J Doe"
$ws.Range("J2").Value = 0
$ws.Range("K4").Value = "Grand Canyon NP"
$ws.Range("I4").Value = "J Doe"
$ws.Range("A4").Value = "GRCA120730"

# Move the selection/active cell as recorded in the saved sheet view.
[void]$ws.Range("L11").Select()
